# Fruta / hortaliza, semanal
# The source data for this sheet was re-pulled and the detail rows (2-48)
# ended up in a different order than before. Column A:T content per record
# is unchanged -- only the row each record sits on changes. Re-order the
# rows to match the refreshed extract by permuting row 2..48 using the
# mapping: new row N gets the content that used to live on row Map[N].

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow  = 48
$lastCol      = 20   # column T

# Map[i] = the ORIGINAL row number whose data should end up at
# data-row index i (i = 0 corresponds to worksheet row 2, etc.)
$map = @(48,46,19,20,18,4,5,37,38,23,26,6,7,33,34,35,47,27,28,39,40,12,10,16,17,8,9,30,31,24,21,22,25,42,43,44,3,13,45,2,11,36,32,29,14,15,41)

# Snapshot all current values (rows 2..48, columns 1..20) before writing
# anything back, so source rows are not clobbered before they are read.
$rowCount = $lastDataRow - $firstDataRow + 1
$snapshot = @{}
for ($i = 0; $i -lt $rowCount; $i++) {
    $srcRow = $firstDataRow + $i
    $rowVals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals += ,$ws.Cells.Item($srcRow, $c).Value2
    }
    $snapshot[$srcRow] = $rowVals
}

# Write the snapshot back out in the permuted order.
for ($i = 0; $i -lt $rowCount; $i++) {
    $destRow = $firstDataRow + $i
    $srcRow  = $map[$i]
    $rowVals = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $rowVals[$c - 1]
    }
}
